$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leaderboard update Cycle 1 - append 5 new result rows (139-143)

$rows = @(
    @{ Row=139; A="Staryy";    C="Staryy";    D=0.28055555555555556; E="Tartaglia "; F="Bennett";   G="Xiangling"; H="Kazuha";  I="Itto";     J="Gorou";  K="TravelerGeo";      L="Kuki" },
    @{ Row=140; A="Not Logic"; C="Not Logic"; D=0.18333333333333335; E="Nilou";      F="Nahida";    G="Kokomi";    H="Collei";  I="HuTao";    J="Xingqiu";K="Mona";             L="Zhongli" },
    @{ Row=141; A="Kurogami";  C="Kurogami";  D=0.15069444444444444; E="Raiden";     F="Kokomi";    G="Collei";    H="Sucrose"; I="Eula";     J="Xinyan"; K="TravelerElectro";  L="Eula" },
    @{ Row=142; A="Yangi";     C="Yangi";     D=0.15138888888888888; E="Nilou";      F="Xingqiu";   G="Yaoyao";    H="Collei";  I="Baizhu";   J="Yae";    K="Fischl";           L="Sucrose" },
    @{ Row=143; A="Ghosted";   C="Ghosted";   D=0.24027777777777778; E="Ayato";      F="Rosaria";   G="Ganyu";     H="Jean";    I="Alhaitham";J="Fischl"; K="YunJin";           L="Zhongli" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = 1
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("D$n").NumberFormat = "h:mm"
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H
    $ws.Range("I$n").Value = $r.I
    $ws.Range("J$n").Value = $r.J
    $ws.Range("K$n").Value = $r.K
    $ws.Range("L$n").Value = $r.L
}

# Match the author's final view state (scroll position + active selection)
$aw = $ws.Application.ActiveWindow
$aw.ScrollRow = 127
$aw.ScrollColumn = 2
$ws.Range("K143").Select()
